# Add the three missing Facebook community rows to the "communities" sheet
# and refresh the JSON-snippet helper formula in column E to cover the new
# rows, matching the upstream commit "Added missing Facebook data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("communities")

# --- New row 7: Facebook "Chat With a Street Epistemologist" Group -------
$ws.Range("A7").Value = "logo_facebook"
$ws.Range("B7").Value = 'Facebook "Chat With a Street Epistemologist" Group'
$ws.Range("C7").Value = "This group is for facilitating SE discussions over particular beliefs."

# --- New row 8: Facebook Critique SE Group --------------------------------
$ws.Range("A8").Value = "logo_facebook"
$ws.Range("B8").Value = "Facebook Critique SE Group"

# --- New row 9: Public SE Facebook Page -----------------------------------
$ws.Range("A9").Value = "logo_facebook"
$ws.Range("B9").Value = "Public SE Facebook Page"

# Descriptions for rows 8 & 9
$ws.Range("C8").Value = "This group is for raising and discussing any critiques of SE so that we can further learn about and improve the method."
$ws.Range("C9").Value = "A public page for believers and non-believers to discuss topics related to Street Epistemology."

# URLs (row 8 reuses the same group URL as row 7 in the source data)
$ws.Range("D7").Value = "https://www.facebook.com/groups/ChatWithAStreetEpistemologist/"
$ws.Range("D8").Value = "https://www.facebook.com/groups/ChatWithAStreetEpistemologist/"
$ws.Range("D9").Value = "https://www.facebook.com/StreetEpistemology"

# Extend the helper JSON-builder formula from E2:E6 down through E2:E9
$formula = '="{
    logo: """&A2&""",
    title: """&B2&""",
    description: """&SUBSTITUTE(C2,"""","\""")&""",
    url: """&D2&"""
},"'
$ws.Range("E2:E9").Formula = $formula

# Re-setting the multi-line formula marks rows 1:9 as "custom height" (Excel's
# auto-wrap side effect) even though the source file keeps the default row
# height; AutoFit puts every row back to the un-flagged default.
$ws.Range("A1:E9").EntireRow.AutoFit()

# Match the saved selection/active cell from the target workbook
$ws.Activate()
$ws.Range("J10").Select()
